# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# The "Metadata" sheet lists FHIR ValueSet properties as Property/Value pairs.
# The "Name" row (A4) was missing its Value (B4) -- fill it in with the
# resource's machine name. Also refresh the "Date" row (A8/B8) value to the
# regenerated timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Name -> Value (row 4) was blank; set it to the resource name.
$ws.Range("B4").Value = "CompetenceexclusiveVs"

# Date -> Value (row 8) bumped to the new generation timestamp.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
